# Apply cryptos list update (Mon May 27 02:57:44 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number/date and must be forced to text
# so Excel does not silently convert "7.50" -> 7.5, "0.0000166" -> 1.66E-5, etc.
$textForceCells = @(
    @{ Cell = "D5"; Value = "605.35" },
    @{ Cell = "D6"; Value = "164.29" },
    @{ Cell = "D9"; Value = "0.532" },
    @{ Cell = "D11"; Value = "6.38" },
    @{ Cell = "D13"; Value = "37.13" },
    @{ Cell = "D18"; Value = "7.50" },
    @{ Cell = "D20"; Value = "17.15" },
    @{ Cell = "D21"; Value = "11.21" },
    @{ Cell = "D22"; Value = "487.85" },
    @{ Cell = "D24"; Value = "0.0000166" },
    @{ Cell = "D25"; Value = "84.46" },
    @{ Cell = "D26"; Value = "2.26" },
    @{ Cell = "D28"; Value = "10.11" },
    @{ Cell = "D32"; Value = "7.89" },
    @{ Cell = "D33"; Value = "32.47" },
    @{ Cell = "D34"; Value = "2.38" },
    @{ Cell = "D40"; Value = "0.999" },
    @{ Cell = "D42"; Value = "3.04" },
    @{ Cell = "D43"; Value = "439.68" },
    @{ Cell = "D45"; Value = "48.46" },
    @{ Cell = "D47"; Value = "8.47" },
    @{ Cell = "D48"; Value = "27.81" },
    @{ Cell = "D50"; Value = "0.0358" },
    @{ Cell = "D51"; Value = "141.17" }
)

foreach ($item in $textForceCells) {
    $rng = $ws.Range($item.Cell)
    $rng.NumberFormat = "@"
    $rng.Value = $item.Value
    $rng.Style = "Normal"
}

# Remaining cells (text / links / percentages / already-non-numeric strings)
# can be assigned directly.
$ws.Range("D2").Value = "69.126.48"
$ws.Range("E2").Value = "  +0.28%  "
$ws.Range("D3").Value = "3.920.70"
$ws.Range("E3").Value = "  +4.83%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("E5").Value = "  +0.64%  "
$ws.Range("E6").Value = "  -0.99%  "
$ws.Range("D7").Value = "3.915.86"
$ws.Range("E7").Value = "  +4.77%  "
$ws.Range("E8").Value = "  +0.09%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -1.58%  "
$ws.Range("E11").Value = "  -0.60%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("E14").Value = "  -0.90%  "
$ws.Range("D15").Value = "4.574.20"
$ws.Range("E15").Value = "  +4.76%  "
$ws.Range("D16").Value = "3.900.35"
$ws.Range("E16").Value = "  +4.14%  "
$ws.Range("D17").Value = "69.189.01"
$ws.Range("E17").Value = "  +0.25%  "
$ws.Range("E18").Value = "  +0.82%  "
$ws.Range("E19").Value = "  -0.41%  "
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  +0.62%  "
$ws.Range("E22").Value = "  -0.78%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  +11.94%  "
$ws.Range("E25").Value = "  -0.17%  "
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -1.52%  "
$ws.Range("E28").Value = "  +0.45%  "
$ws.Range("E29").Value = "  +0.04%  "
$ws.Range("E30").Value = "  -0.46%  "
$ws.Range("D31").Value = "4.074.08"
$ws.Range("E31").Value = "  +4.85%  "
$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("D35").Value = "3.865.12"
$ws.Range("E35").Value = "  +5.15%  "
$ws.Range("E36").Value = "  -0.47%  "
$ws.Range("E37").Value = "  +3.25%  "
$ws.Range("E38").Value = "  +1.27%  "
$ws.Range("E39").Value = "  +0.18%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  -1.25%  "
$ws.Range("E42").Value = "  -4.77%  "
$ws.Range("E43").Value = "  +3.20%  "
$ws.Range("E44").Value = "  +0.24%  "
$ws.Range("E45").Value = "  -0.26%  "
$ws.Range("E46").Value = "  -0.02%  "
$ws.Range("B47").Value = "Cosmos"
$ws.Range("C47").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("E47").Value = "  +0.23%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("E48").Value = "  +18.32%  "
$ws.Range("D49").Value = "2.850.96"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("E50").Value = "  +2.07%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E51").Value = "  +0.04%  "
